$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 943.7778
$ws.Range("I12").Value = 249.66667
$ws.Range("J12").Value = 2332
$ws.Range("K12").Value = 249.66667
$ws.Range("L12").Value = 2332
$ws.Range("M12").Value = -79.66667000000001
$ws.Range("N12").Value = -2672
$ws.Range("H40").Value = 4693.316
$ws.Range("J40").Value = 7124.75
$ws.Range("L40").Value = 7124.75
$ws.Range("N40").Value = -7474.75
$ws.Range("H41").Value = 544.5
$ws.Range("I41").Value = 544.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 544.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -104.5
$ws.Range("N41").ClearContents()
$ws.Range("H100").Value = 1238.3478
$ws.Range("I100").Value = 975.125
$ws.Range("K100").Value = 975.125
$ws.Range("M100").Value = -434.125
$ws.Range("H131").Value = 145536.72
$ws.Range("I131").Value = 145536.72
$ws.Range("K131").Value = 436610.16
$ws.Range("M131").Value = -431570.16
$ws.Range("H133").Value = 77272.25
$ws.Range("J133").Value = 77272.25
$ws.Range("L133").Value = 77272.25
$ws.Range("N133").Value = -87392.25
$ws.Range("H135").Value = 589.4737
$ws.Range("I135").Value = 589.4737
$ws.Range("K135").Value = 5305.263300000001
$ws.Range("M135").Value = -2770.263300000001
$ws.Range("H138").Value = 1857.1364
$ws.Range("J138").Value = 2525.3667
$ws.Range("L138").Value = 7576.1001
$ws.Range("N138").Value = -17856.1001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3911.3914
$ws.Range("I32").Value = 3713.372
$ws.Range("K32").Value = 3713.372
$ws.Range("M32").Value = -3426.372
$ws.Range("H61").Value = 6150.452
$ws.Range("I61").Value = 5164.4243
$ws.Range("K61").Value = 5164.4243
$ws.Range("M61").Value = -4952.4243
$ws.Range("H102").Value = 7333.222
$ws.Range("I102").Value = 4999.8335
$ws.Range("K102").Value = 4999.8335
$ws.Range("M102").Value = -3377.8335
$ws.Range("H132").Value = 7571.357
$ws.Range("I132").Value = 6463.5654
$ws.Range("K132").Value = 19390.6962
$ws.Range("M132").Value = -16860.6962
$ws.Range("H136").Value = 6150.452
$ws.Range("I136").Value = 5164.4243
$ws.Range("K136").Value = 15493.2729
$ws.Range("M136").Value = -12943.2729

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5650.222
$ws.Range("I99").Value = 4520.4
$ws.Range("K99").Value = 4520.4
$ws.Range("M99").Value = -3022.4
$ws.Range("H107").Value = 3212.7144
$ws.Range("I107").Value = 3212.7144
$ws.Range("K107").Value = 3212.7144
$ws.Range("M107").Value = -1292.7144
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 363
$ws.Range("I22").Value = 353.83334
$ws.Range("K22").Value = 353.83334
$ws.Range("M22").Value = -3.833340000000021
$ws.Range("H31").Value = 4622.75
$ws.Range("I31").Value = 3728.3333
$ws.Range("J31").Value = 5654.769
$ws.Range("K31").Value = 3728.3333
$ws.Range("L31").Value = 5654.769
$ws.Range("M31").Value = -3433.3333
$ws.Range("N31").Value = -6244.769
$ws.Range("H34").Value = 4622.75
$ws.Range("I34").Value = 3728.3333
$ws.Range("J34").Value = 5654.769
$ws.Range("K34").Value = 3728.3333
$ws.Range("L34").Value = 5654.769
$ws.Range("M34").Value = -3526.3333
$ws.Range("N34").Value = -6058.769
$ws.Range("H58").Value = 5287.4688
$ws.Range("I58").Value = 3670.8333
$ws.Range("K58").Value = 3670.8333
$ws.Range("M58").Value = -3467.8333
$ws.Range("H136").Value = 5287.4688
$ws.Range("I136").Value = 3670.8333
$ws.Range("K136").Value = 11012.4999
$ws.Range("M136").Value = -8462.499899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 31500276
$ws.Range("I11").Value = 42000030
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 126000090
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -125999950
$ws.Range("N11").Value = -3280
$ws.Range("H52").Value = 2222
$ws.Range("J52").Value = 2222
$ws.Range("L52").Value = 6666
$ws.Range("N52").Value = -7198
$ws.Range("H131").Value = 29414188
$ws.Range("I131").Value = 100000910
$ws.Range("J131").Value = 3054.1667
$ws.Range("K131").Value = 300002730
$ws.Range("L131").Value = 9162.500100000001
$ws.Range("M131").Value = -299997690
$ws.Range("N131").Value = -19242.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H51").Value = 85000
$ws.Range("J51").Value = 85000
$ws.Range("L51").Value = 85000
$ws.Range("N51").Value = -86018
$ws.Range("H70").Value = 14144.533
$ws.Range("I70").Value = 11579.6
$ws.Range("J70").Value = 15427
$ws.Range("K70").Value = 11579.6
$ws.Range("L70").Value = 15427
$ws.Range("M70").Value = -11309.6
$ws.Range("N70").Value = -15967
$ws.Range("H73").Value = 14144.533
$ws.Range("I73").Value = 11579.6
$ws.Range("J73").Value = 15427
$ws.Range("K73").Value = 11579.6
$ws.Range("L73").Value = 15427
$ws.Range("M73").Value = -10643.6
$ws.Range("N73").Value = -17299
$ws.Range("H132").Value = 2167.6
$ws.Range("I132").Value = 1973.5217
$ws.Range("J132").Value = 4399.5
$ws.Range("K132").Value = 5920.5651
$ws.Range("L132").Value = 13198.5
$ws.Range("M132").Value = -3390.5651
$ws.Range("N132").Value = -18258.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1539.8572
$ws.Range("I22").Value = 1595.8
$ws.Range("K22").Value = 1595.8
$ws.Range("M22").Value = -1300.8
$ws.Range("H27").Value = 1539.8572
$ws.Range("I27").Value = 1595.8
$ws.Range("K27").Value = 1595.8
$ws.Range("M27").Value = -1488.8
$ws.Range("H55").Value = 624.53845
$ws.Range("I55").Value = 675.7
$ws.Range("J55").Value = 454
$ws.Range("K55").Value = 675.7
$ws.Range("L55").Value = 454
$ws.Range("M55").Value = -502.7
$ws.Range("N55").Value = -800
$ws.Range("H68").Value = 1998.1428
$ws.Range("I68").Value = 1996
$ws.Range("K68").Value = 1996
$ws.Range("M68").Value = -1247
$ws.Range("H71").Value = 1998.1428
$ws.Range("I71").Value = 1996
$ws.Range("K71").Value = 9980
$ws.Range("M71").Value = -6236
$ws.Range("H100").Value = 7499.5
$ws.Range("J100").Value = 7499.6665
$ws.Range("L100").Value = 7499.6665
$ws.Range("N100").Value = -8581.666499999999
$ws.Range("H122").Value = 3079.1304
$ws.Range("I122").Value = 2789.1428
$ws.Range("K122").Value = 8367.428400000001
$ws.Range("M122").Value = -5917.428400000001
$ws.Range("H132").Value = 7976.9556
$ws.Range("I132").Value = 7749.026
$ws.Range("K132").Value = 23247.078
$ws.Range("M132").Value = -20717.078
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2229.8333
$ws.Range("I81").Value = 899.25
$ws.Range("J81").Value = 4891
$ws.Range("K81").Value = 1798.5
$ws.Range("L81").Value = 9782
$ws.Range("M81").Value = -737.5
$ws.Range("N81").Value = -11904
$ws.Range("H84").Value = 2229.8333
$ws.Range("I84").Value = 899.25
$ws.Range("J84").Value = 4891
$ws.Range("K84").Value = 8992.5
$ws.Range("L84").Value = 48910
$ws.Range("M84").Value = -3688.5
$ws.Range("N84").Value = -59518
$ws.Range("H96").Value = 2893.3845
$ws.Range("I96").Value = 2587.8572
$ws.Range("J96").Value = 3249.8333
$ws.Range("K96").Value = 2587.8572
$ws.Range("L96").Value = 3249.8333
$ws.Range("M96").Value = -1214.8572
$ws.Range("N96").Value = -5995.8333
$ws.Range("H107").Value = 999.7143
$ws.Range("I107").Value = 692.4286
$ws.Range("K107").Value = 2077.2858
$ws.Range("M107").Value = -157.2857999999997
$ws.Range("H126").Value = 6195
$ws.Range("J126").Value = 5080.8
$ws.Range("L126").Value = 15242.4
$ws.Range("N126").Value = -20182.4
